$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the existing hyperlink on A49 (it currently holds the URL text)
$ws.Hyperlinks.Delete()

# Current layout (before edit):
#   A47 "Source:"
#   A48 "Le rapport de l'enquête sur les entreprises à Madagascar, INSTAT, p. 9"
#   A49 "http://www.instat.mg/pdf/enq_entreprises_2005.pdf"  (hyperlink style)
#   A50 "" (blank, source style)
#   A53 "INSTAT" (title style)
#   A54 "Institut National de la Statistique de Madagascar (INSTAT), ..." (source style)
#
# Target layout (after edit):
#   A47 "Source:"
#   A48 "" (blank, source style)
#   A49 "Le rapport de l'enquête sur les entreprises à Madagascar, INSTAT, p. 9" (source style)
#   A50 "" (blank, source style)  -- unchanged
#   A51 "http://www.instat.mg/pdf/enq_entreprises_2005.pdf" (source style, plain text, no hyperlink)
#   A54 "INSTAT" (title style)
#   A55 "INSTAT" (source style)

# Move A48's current text down to A49 first (so we don't lose it while editing A48)
$ws.Range("A49").Value = $ws.Range("A48").Value
$ws.Range("A49").Style = "source"

# Clear A48 and make it blank with source style
$ws.Range("A48").Value = ""
$ws.Range("A48").Style = "source"

# New row A51 holds the plain-text URL (no hyperlink, source style)
$ws.Range("A51").Value = "http://www.instat.mg/pdf/enq_entreprises_2005.pdf"
$ws.Range("A51").Style = "source"

# Move "INSTAT" title from A53 to A54
$ws.Range("A54").Value = $ws.Range("A53").Value
$ws.Range("A54").Style = "title"
$ws.Range("A53").Value = ""

# A55 becomes "INSTAT" with source style (replaces old long citation previously at A54)
$ws.Range("A55").Value = "INSTAT"
$ws.Range("A55").Style = "source"
